$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workflow")

# Eject variable: replace the single "employee_name" eject-variable expression used on
# both the GetEmployeeName (row 2) and GetLastName (row 3) rows with a multi-field
# eject-variable expression on row 2 only, and clear it from row 3.
$ws.Range("G2").Value = "FirstName=data[0].first_name,LastName=data[0].last_name"
$ws.Range("G3").ClearContents()

# Widen column G to fit the longer eject-variable text.
$ws.Columns.Item(7).ColumnWidth = 37.67

# Update the active selection on the sheet.
$ws.Activate()
$ws.Range("G4").Select()
